$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 386.9
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H9").Value = 946.6667
$ws.Range("J9").Value = 190
$ws.Range("L9").Value = 190
$ws.Range("N9").Value = -528
$ws.Range("H28").Value = 2291.7
$ws.Range("I28").Value = 1879.7778
$ws.Range("K28").Value = 1879.7778
$ws.Range("M28").Value = -1394.7778
$ws.Range("H38").Value = 865.13336
$ws.Range("I38").Value = 784.0714
$ws.Range("J38").Value = 2000
$ws.Range("K38").Value = 2352.2142
$ws.Range("L38").Value = 6000
$ws.Range("M38").Value = -1980.2142
$ws.Range("N38").Value = -6744
$ws.Range("H43").Value = 1849
$ws.Range("I43").Value = 999
$ws.Range("K43").Value = 999
$ws.Range("M43").Value = -930
$ws.Range("H69").Value = 16849.9
$ws.Range("I69").Value = 14499.667
$ws.Range("J69").Value = 17857.143
$ws.Range("K69").Value = 43499.001
$ws.Range("L69").Value = 53571.429
$ws.Range("M69").Value = -42625.001
$ws.Range("N69").Value = -55319.429
$ws.Range("H72").Value = 16849.9
$ws.Range("I72").Value = 14499.667
$ws.Range("J72").Value = 17857.143
$ws.Range("K72").Value = 130497.003
$ws.Range("L72").Value = 160714.287
$ws.Range("M72").Value = -126129.003
$ws.Range("N72").Value = -169450.287
$ws.Range("H138").Value = 2245.606
$ws.Range("I138").Value = 2028.8462
$ws.Range("J138").Value = 2386.5
$ws.Range("K138").Value = 6086.5386
$ws.Range("L138").Value = 7159.5
$ws.Range("M138").Value = -946.5385999999999
$ws.Range("N138").Value = -17439.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2875.7727
$ws.Range("I45").Value = 2681.7896
$ws.Range("J45").Value = 4104.3335
$ws.Range("K45").Value = 2681.7896
$ws.Range("L45").Value = 4104.3335
$ws.Range("M45").Value = -2304.7896
$ws.Range("N45").Value = -4858.3335
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H97").Value = 2880.7273
$ws.Range("I97").Value = 1332
$ws.Range("K97").Value = 1332
$ws.Range("M97").Value = -836
$ws.Range("H110").Value = 17999.143
$ws.Range("I110").Value = 13499
$ws.Range("K110").Value = 13499
$ws.Range("M110").Value = -11454

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3724
$ws.Range("I94").Value = 3099
$ws.Range("K94").Value = 3099
$ws.Range("M94").Value = -2648

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 103.35
$ws.Range("I7").Value = 61.42105
$ws.Range("J7").Value = 900
$ws.Range("K7").Value = 61.42105
$ws.Range("L7").Value = 900
$ws.Range("M7").Value = 51.57895
$ws.Range("N7").Value = -1126
$ws.Range("H26").Value = 20996
$ws.Range("J26").Value = 20996
$ws.Range("L26").Value = 20996
$ws.Range("N26").Value = -21570
$ws.Range("H31").Value = 9224.362999999999
$ws.Range("I31").Value = 12868.714
$ws.Range("J31").Value = 2846.75
$ws.Range("K31").Value = 12868.714
$ws.Range("L31").Value = 2846.75
$ws.Range("M31").Value = -12573.714
$ws.Range("N31").Value = -3436.75
$ws.Range("H34").Value = 9224.362999999999
$ws.Range("I34").Value = 12868.714
$ws.Range("J34").Value = 2846.75
$ws.Range("K34").Value = 12868.714
$ws.Range("L34").Value = 2846.75
$ws.Range("M34").Value = -12666.714
$ws.Range("N34").Value = -3250.75
$ws.Range("H41").Value = 8450
$ws.Range("I41").Value = 8450
$ws.Range("K41").Value = 8450
$ws.Range("M41").Value = -8022
$ws.Range("H58").Value = 3023.4
$ws.Range("I58").Value = 2006.2858
$ws.Range("J58").Value = 5396.6665
$ws.Range("K58").Value = 2006.2858
$ws.Range("L58").Value = 5396.6665
$ws.Range("M58").Value = -1803.2858
$ws.Range("N58").Value = -5802.6665
$ws.Range("H86").Value = 14776
$ws.Range("I86").Value = 8997.857
$ws.Range("K86").Value = 8997.857
$ws.Range("M86").Value = -7874.857
$ws.Range("H89").Value = 14776
$ws.Range("I89").Value = 8997.857
$ws.Range("K89").Value = 44989.285
$ws.Range("M89").Value = -39373.285
$ws.Range("H136").Value = 3023.4
$ws.Range("I136").Value = 2006.2858
$ws.Range("J136").Value = 5396.6665
$ws.Range("K136").Value = 6018.857400000001
$ws.Range("L136").Value = 16189.9995
$ws.Range("M136").Value = -3468.857400000001
$ws.Range("N136").Value = -21289.9995

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2252.375
$ws.Range("I5").Value = 2261.3076
$ws.Range("K5").Value = 6783.9228
$ws.Range("M5").Value = -6671.9228
$ws.Range("H40").Value = 137.75
$ws.Range("I40").Value = 33
$ws.Range("J40").Value = 242.5
$ws.Range("K40").Value = 132
$ws.Range("L40").Value = 970
$ws.Range("M40").Value = -63
$ws.Range("N40").Value = -1108
$ws.Range("H68").Value = 630
$ws.Range("I68").Value = 630
$ws.Range("K68").Value = 1890
$ws.Range("M68").Value = -1079
$ws.Range("H71").Value = 630
$ws.Range("I71").Value = 630
$ws.Range("K71").Value = 5670
$ws.Range("M71").Value = -1614
$ws.Range("H93").Value = 9055.875
$ws.Range("I93").Value = 861.75
$ws.Range("K93").Value = 2585.25
$ws.Range("M93").Value = -713.25
$ws.Range("H135").Value = 2252.375
$ws.Range("I135").Value = 2261.3076
$ws.Range("K135").Value = 20351.7684
$ws.Range("M135").Value = -17816.7684

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 27130.5
$ws.Range("J39").Value = 24261
$ws.Range("L39").Value = 24261
$ws.Range("N39").Value = -25325
$ws.Range("H113").Value = 5833.3335
$ws.Range("I113").Value = 3000
$ws.Range("K113").Value = 3000
$ws.Range("M113").Value = -830

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 892.9167
$ws.Range("I22").Value = 843.125
$ws.Range("J22").Value = 992.5
$ws.Range("K22").Value = 843.125
$ws.Range("L22").Value = 992.5
$ws.Range("M22").Value = -548.125
$ws.Range("N22").Value = -1582.5
$ws.Range("H27").Value = 892.9167
$ws.Range("I27").Value = 843.125
$ws.Range("J27").Value = 992.5
$ws.Range("K27").Value = 843.125
$ws.Range("L27").Value = 992.5
$ws.Range("M27").Value = -736.125
$ws.Range("N27").Value = -1206.5
$ws.Range("H61").Value = 1100
$ws.Range("I61").Value = 1100
$ws.Range("K61").Value = 1100
$ws.Range("M61").Value = -898
$ws.Range("H82").Value = 13300
$ws.Range("I82").Value = 13300
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 13300
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -12939
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 13300
$ws.Range("I85").Value = 13300
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 13300
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -12052
$ws.Range("N85").ClearContents()
$ws.Range("H113").Value = 1100
$ws.Range("I113").Value = 1100
$ws.Range("K113").Value = 1100
$ws.Range("M113").Value = 1070
$ws.Range("H132").Value = 4943.3335
$ws.Range("I132").Value = 4641.4287
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 13924.2861
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -11394.2861
$ws.Range("N132").Value = -23060

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2300.1667
$ws.Range("I81").Value = 1759.8
$ws.Range("J81").Value = 5002
$ws.Range("K81").Value = 3519.6
$ws.Range("L81").Value = 10004
$ws.Range("M81").Value = -2458.6
$ws.Range("N81").Value = -12126
$ws.Range("H84").Value = 2300.1667
$ws.Range("I84").Value = 1759.8
$ws.Range("J84").Value = 5002
$ws.Range("K84").Value = 17598
$ws.Range("L84").Value = 50020
$ws.Range("M84").Value = -12294
$ws.Range("N84").Value = -60628
$ws.Range("H109").Value = 70000
$ws.Range("J109").Value = 70000
$ws.Range("L109").Value = 70000
$ws.Range("N109").Value = -72774
$ws.Range("H124").Value = 23975.4
$ws.Range("J124").Value = 23975.4
$ws.Range("L124").Value = 23975.4
$ws.Range("N124").Value = -33795.4
$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -49840
